$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns E..T correspond to rows in the array below, in order:
# E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T
$data = @(
    @(2, 3, 1, 7.572976999999999, 22.718931, 0.4497670593913077, 0.4497670593913078, 3, 1, 261.380203, 784.1406089999999, 0.6968677182772199, 0.6968677182772199, 1979.426265574331, 17814.83639016898, 0.3134281444342755, 0.3134281444342755),
    @(3, 3, 1, 7.572976999999999, 22.718931, 0.4497670593913077, 0.4497670593913078, 3, 1, 31.999428, 95.998284, 0.08531391482826334, 0.08531391482826335, 242.330932257156, 2180.978390314404, 0.03837138859746848, 0.0383713885974685),
    @(4, 3, 1, 7.572976999999999, 22.718931, 0.4497670593913077, 0.4497670593913078, 3, 1, 81.69901900000001, 245.097057, 0.2178183668945166, 0.2178183668945167, 618.7047918095629, 5568.343126286067, 0.09796752635956373, 0.09796752635956375),
    @(5, 3, 1, 2.784025666666667, 8.352077, 0.1653462089435359, 0.1653462089435359, 3, 1, 261.380203, 784.1406089999999, 0.6968677182772199, 0.6968677182772199, 727.6891939105436, 6549.202745194892, 0.1152244353522704, 0.1152244353522704),
    @(6, 3, 1, 2.784025666666667, 8.352077, 0.1653462089435359, 0.1653462089435359, 3, 1, 31.999428, 95.998284, 0.08531391482826334, 0.08531391482826335, 89.08722887065198, 801.785059835868, 0.01410633238698506, 0.01410633238698506),
    @(7, 3, 1, 2.784025666666667, 8.352077, 0.1653462089435359, 0.1653462089435359, 3, 1, 81.69901900000001, 245.097057, 0.2178183668945166, 0.2178183668945167, 227.4521658374877, 2047.069492537389, 0.03601544120428052, 0.03601544120428053),
    @(8, 3, 1, 6.480551000000001, 19.441653, 0.3848867316651562, 0.3848867316651562, 3, 1, 261.380203, 784.1406089999999, 0.6968677182772199, 0.6968677182772199, 1693.887735931853, 15244.98962338668, 0.268215138490674, 0.2682151384906741),
    @(9, 3, 1, 6.480551000000001, 19.441653, 0.3848867316651562, 0.3848867316651562, 3, 1, 31.999428, 95.998284, 0.08531391482826334, 0.08531391482826335, 207.373925124828, 1866.365326123452, 0.03283619384380978, 0.03283619384380979),
    @(10, 3, 1, 6.480551000000001, 19.441653, 0.3848867316651562, 0.3848867316651562, 3, 1, 81.69901900000001, 245.097057, 0.2178183668945166, 0.2178183668945167, 529.4546592794692, 4765.091933515222, 0.08383539933067237, 0.08383539933067238)
)

foreach ($row in $data) {
    $r = $row[0]
    for ($i = 0; $i -lt 16; $i++) {
        $col = $i + 5  # column E = 5
        $ws.Cells.Item($r, $col).Value = $row[$i + 1]
    }
}